$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that contain the SQL query text needing the join-key rename.
$cellRefs = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($ref in $cellRefs) {
    $cell = $ws.Range($ref)
    $text = $cell.Value2
    if ($text -ne $null) {
        $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
        $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
        $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
        $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
        $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
        $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
        $cell.Value2 = $text
    }
}

# Update the view: move the active selection to C7 and scroll row 6 to the top
# (mirrors topLeftCell="A6" / selection activeCell="C7" sqref="C7" in the saved XML).
$ws.Range("C7").Select()
$excel.ActiveWindow.ScrollRow = 6
